$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-All "307 (100.0)" "329 (100.0)"
Replace-All "281 (100.0)" "303 (100.0)"
Replace-All "281 (91.8)" "303 (92.4)"
Replace-All "205 (69.5)" "225 (69.4)"
Replace-All "205 (75.4)" "225 (75.3)"
Replace-All "157 (55.9)" "169 (54.9)"
Replace-All "157 (60.2)" "169 (59.7)"
Replace-All "143 (54.0)" "152 (53.5)"
Replace-All "143 (57.9)" "152 (57.6)"
